$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.491945396372404
$ws.Range("C2").Value = 0.49220524614661
$ws.Range("B3").Value = 0.455865178149482
$ws.Range("C3").Value = 0.454902509849823
$ws.Range("B4").Value = 0.255873391728952
$ws.Range("C4").Value = 0.256717934563359
$ws.Range("B5").Value = 0.630419104437875
$ws.Range("C5").Value = 0.630007896550884
$ws.Range("B6").Value = 0.659485831874964
$ws.Range("C6").Value = 0.660272612487221
$ws.Range("B7").Value = 0.540655018328467
$ws.Range("C7").Value = 0.542690479734615
$ws.Range("B8").Value = 0.737382789402442
$ws.Range("C8").Value = 0.735950890557437
$ws.Range("B9").Value = 0.786324640463071
$ws.Range("C9").Value = 0.786559537344948
$ws.Range("B10").Value = 0.63510138672369
$ws.Range("C10").Value = 0.635498627393652
$ws.Range("B11").Value = 0.734113151723734
$ws.Range("C11").Value = 0.7342148993472
$ws.Range("B12").Value = 0.534678276289036
$ws.Range("C12").Value = 0.534480549863145
$ws.Range("B13").Value = 0.949340103048975
$ws.Range("C13").Value = 0.949378946983765
$ws.Range("B14").Value = 0.56628802854227
$ws.Range("C14").Value = 0.566192776127786
$ws.Range("B15").Value = 0.814871237474843
$ws.Range("C15").Value = 0.815213495873063
$ws.Range("B16").Value = 0.700784933571959
$ws.Range("C16").Value = 0.698654908634361
$ws.Range("B17").Value = 0.635270097218214
$ws.Range("C17").Value = 0.637808022237456
$ws.Range("B18").Value = 0.836909470811762
$ws.Range("C18").Value = 0.839226701597952
$ws.Range("B19").Value = 0.645203803383275
$ws.Range("C19").Value = 0.645328115994307
$ws.Range("B20").Value = 0.425471624783544
$ws.Range("C20").Value = 0.4254017981675
$ws.Range("B21").Value = 0.460626918026747
$ws.Range("C21").Value = 0.460150786390864
$ws.Range("B22").Value = 0.51612468554495
$ws.Range("C22").Value = 0.517378959828486
$ws.Range("B23").Value = 0.705839161308343
$ws.Range("C23").Value = 0.704527710256888
$ws.Range("B24").Value = 0.587123294766919
$ws.Range("C24").Value = 0.588589400461793
$ws.Range("B25").Value = 0.416354043259093
$ws.Range("C25").Value = 0.419346820527562
$ws.Range("B26").Value = 0.863622233182104
$ws.Range("C26").Value = 0.867860738206216
$ws.Range("B27").Value = 0.637928451019375
$ws.Range("C27").Value = 0.622679680185949
$ws.Range("B28").Value = 0.862216976188318
$ws.Range("C28").Value = 0.874300338159811
$ws.Range("B29").Value = 0.855080362371329
$ws.Range("C29").Value = 0.852445616348055
$ws.Range("B30").Value = 0.956882332398853
$ws.Range("C30").Value = 0.94517020180737
$ws.Range("B31").Value = 0.962942446106548
$ws.Range("C31").Value = 0.961869939932097
$ws.Range("B32").Value = 0.856609675674819
$ws.Range("C32").Value = 0.848045421414048
$ws.Range("B33").Value = 0.980625221460318
$ws.Range("C33").Value = 0.981413699672006
$ws.Range("B34").Value = 0.956199985609417
$ws.Range("C34").Value = 0.956417897829349
